$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.267.56'
$ws.Range("E2").Value = '  +0.32%  '
$ws.Range("D3").Value = '1.851.39'
$ws.Range("E3").Value = '  +0.98%  '
$ws.Range("D4").Value = "'0.9994"
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = "'241.41"
$ws.Range("E5").Value = '  +0.45%  '
$ws.Range("D6").Value = "'0.6738"
$ws.Range("E6").Value = '  -1.63%  '
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = "'0.07449"
$ws.Range("E8").Value = '  -0.27%  '
$ws.Range("D9").Value = "'0.2971"
$ws.Range("E9").Value = '  -1.48%  '
$ws.Range("D10").Value = "'22.96"
$ws.Range("E10").Value = '  -0.72%  '
$ws.Range("D11").Value = "'0.07735"
$ws.Range("E11").Value = '  +0.99%  '
$ws.Range("D12").Value = '1.834.62'
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").Value = "'5.027"
$ws.Range("E13").Value = '  -0.66%  '
$ws.Range("D14").Value = "'0.6804"
$ws.Range("E14").Value = '  -0.25%  '
$ws.Range("D15").Value = "'86.42"
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("D16").Value = "'6.186"
$ws.Range("E16").Value = '  -0.33%  '
$ws.Range("D17").Value = '29.222.07'
$ws.Range("E17").Value = '  +0.21%  '
$ws.Range("D18").Value = "'0.000008314"
$ws.Range("E18").Value = '  +1.58%  '
$ws.Range("D19").Value = "'229.42"
$ws.Range("E19").Value = '  +1.30%  '
$ws.Range("D20").Value = "'12.59"
$ws.Range("E20").Value = '  +0.42%  '
$ws.Range("D21").Value = "'0.9997"
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("D22").Value = "'7.245"
$ws.Range("E22").Value = '  -2.24%  '
$ws.Range("D23").Value = "'1.000"
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").Value = "'160.96"
$ws.Range("E24").Value = '  +0.83%  '
$ws.Range("D25").Value = "'0.1421"
$ws.Range("E25").Value = '  -2.50%  '
$ws.Range("D26").Value = "'8.719"
$ws.Range("E26").Value = '  -0.50%  '
$ws.Range("D27").Value = "'18.07"
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").Value = "'1.514"
$ws.Range("E28").Value = '  +0.46%  '
$ws.Range("D29").Value = "'4.199"
$ws.Range("E29").Value = '  -1.37%  '
$ws.Range("D30").Value = "'4.090"
$ws.Range("E30").Value = '  -1.26%  '
$ws.Range("D31").Value = "'1.188"
$ws.Range("E31").Value = '  -1.35%  '
$ws.Range("D32").Value = "'0.05335"
$ws.Range("E32").Value = '  +3.58%  '
$ws.Range("D33").Value = "'1.901"
$ws.Range("E33").Value = '  +3.10%  '
$ws.Range("D34").Value = "'0.7588"
$ws.Range("E34").Value = '  -1.10%  '
$ws.Range("D35").Value = "'1.143"
$ws.Range("E35").Value = '  +0.93%  '
$ws.Range("D36").Value = "'2.687"
$ws.Range("E36").Value = '  +0.57%  '
$ws.Range("D37").Value = '1.337.75'
$ws.Range("E37").Value = '  +2.28%  '
$ws.Range("D38").Value = "'0.01807"
$ws.Range("E38").Value = '  -1.58%  '
$ws.Range("D39").Value = "'2.743"
$ws.Range("E39").Value = '  +1.13%  '
$ws.Range("D40").Value = "'0.9273"
$ws.Range("E40").Value = '  -0.81%  '
$ws.Range("D41").Value = "'5.976"
$ws.Range("E41").Value = '  +2.64%  '
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("D43").Value = "'103.72"
$ws.Range("E43").Value = '  -0.56%  '
$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").Value = '1.978.10'
$ws.Range("E44").Value = '  -0.16%  '
$ws.Range("D45").Value = "'0.07829"
$ws.Range("E45").Value = '  +7.37%  '
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").Value = "'0.00000000123"
$ws.Range("E46").Value = '  +2.73%  '
$ws.Range("D47").Value = "'0.5163"
$ws.Range("E47").Value = '  -0.69%  '
$ws.Range("D48").Value = "'1.773"
$ws.Range("E48").Value = '  +0.17%  '
$ws.Range("D49").Value = "'63.94"
$ws.Range("E49").Value = '  -2.32%  '
$ws.Range("D50").Value = "'9.230"
$ws.Range("E50").Value = '  -3.70%  '
$ws.Range("D51").Value = "'0.05945"
$ws.Range("E51").Value = '  +0.50%  '
